# The commit reorders the "Recorded By" value in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where it occurs (cells that already contain only
# "System" or only "dnasr281@gmail.com" are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$ws.Cells.Replace($target, $replacement, [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
